$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 73, pushing the existing rows 73:132 down to 74:133.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new price-report record.
$ws.Range("A73").Value = 10
$ws.Range("B73").Value = "Vega Modelo de Temuco"
$ws.Range("C73").Value = "La Araucanía"
$ws.Range("D73").Value = 44741
$ws.Range("E73").Value = 9
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100104
$ws.Range("H73").Value = "Frutos de pepita"
$ws.Range("I73").Value = 100104001
$ws.Range("J73").Value = "Granada"
$ws.Range("K73").Value = "Wonderfull"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 400
$ws.Range("N73").Value = 13000
$ws.Range("O73").Value = 13000
$ws.Range("P73").Value = 13000
$ws.Range("Q73").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R73").Value = "Región de O'Higgins"
$ws.Range("S73").Value = 1300
$ws.Range("T73").Value = 10
